# Applies the "Now you can add and remove many graphs" edit:
#  - Slide 3, body placeholder, paragraph 1: expands the short
#    "Расширение для Microsoft Visio - VSTO (C#)" bullet into a longer
#    sentence, split across several runs (mixed ru / en-US language runs,
#    matching how PowerPoint's autocorrect/spellcheck splits typed text).
#  - Slide 3, body placeholder, paragraph 2: expands the short
#    "Чтение DOT файла - Graphviz4Net + ANTLR" bullet into a longer
#    sentence (two runs).
#  - Slide 5, body placeholder, paragraph 1: expands
#    "Краткосрочная перспектива:" into "Краткосрочная перспектива (до
#    защиты):" (two runs).
#  - Slide 5, body placeholder, last paragraph: merges the trailing
#    ", чтоб не жаловался ... не " + "объявленных вершин" runs into one.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 ("Как я это делал") - body placeholder is shape 2.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange

# --- Paragraph 1: "Расширение для Microsoft Visio - VSTO (C#)" ---
$para1 = $body3.Paragraphs(1, 1)

# Rebuild the paragraph text/runs from the end backwards: set the whole
# paragraph text to the last chunk, then repeatedly InsertBefore() each
# earlier chunk and (re)stamp the paragraph's LanguageID right after each
# insert. In this COM host, TextRange.LanguageID always stamps the FIRST
# run of the owning paragraph, so doing it immediately after prepending a
# chunk (so it's temporarily the first/only new run) lands the language
# on the correct run.

$para1.Text = "(C#)"
$para1.LanguageID = "ru"

$para1.InsertBefore(" ")
$para1.LanguageID = "ru"

$para1.InsertBefore("ffice")
$para1.LanguageID = "en-US"

$para1.InsertBefore("O")
$para1.LanguageID = "ru"

$para1.InsertBefore(" for ")
$para1.LanguageID = "en-US"

$para1.InsertBefore("ools")
$para1.LanguageID = "en-US"

$para1.InsertBefore("T")
$para1.LanguageID = "ru"

$para1.InsertBefore(" ")
$para1.LanguageID = "en-US"

$para1.InsertBefore("tudio")
$para1.LanguageID = "en-US"

$para1.InsertBefore("S")
$para1.LanguageID = "ru"

$para1.InsertBefore(" ")
$para1.LanguageID = "en-US"

$para1.InsertBefore("isual")
$para1.LanguageID = "en-US"

$para1.InsertBefore("написано при помощи инструмента V")
$para1.LanguageID = "ru"

$para1.InsertBefore("Расширение для Microsoft Visio ")
$para1.LanguageID = "ru"

# --- Paragraph 2: "Чтение DOT файла - Graphviz4Net + ANTLR" ---
$para2 = $body3.Paragraphs(2, 1)
$prefixLen = "Чтение DOT файла ".Length
$tailLen = $para2.Text.Length - $prefixLen
$tail2 = $para2.Characters($prefixLen + 1, $tailLen)
$tail2.Text = "осуществляется при помощи библиотеки Graphviz4Net, включающей в себя парсер ANTLR"

# ---------------------------------------------------------------------
# Slide 5 ("Что дальше?") - body placeholder is shape 2.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange

# --- Paragraph 1: "Краткосрочная перспектива:" ---
$para1b = $body5.Paragraphs(1, 1)
$prefixLen2 = "Краткосрочная ".Length
$tailLen2 = $para1b.Text.Length - $prefixLen2
$tail1b = $para1b.Characters($prefixLen2 + 1, $tailLen2)
$tail1b.Text = "перспектива (до защиты):"

# --- Last paragraph: merge ", чтоб ... не " + "объявленных вершин" ---
$lastParaIdx = $body5.Paragraphs().Count
$paraLast = $body5.Paragraphs($lastParaIdx, 1)
$mergedSuffix = ", чтоб не жаловался на использование ранее не объявленных вершин"
$mergeStart = $paraLast.Text.Length - $mergedSuffix.Length + 1
$mergeRange = $paraLast.Characters($mergeStart, $mergedSuffix.Length)
$mergeRange.Text = $mergedSuffix
